$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added at the top of the "Camote" table
# (row 95). Every existing record from row 95 down to row 106 shifts down by
# one row (row 106's old data ends up at row 107); Excel's own row-insert
# semantics handle that shift for us.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new observation.
$ws.Cells.Item(95, 1).Value = 10
$ws.Cells.Item(95, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(95, 3).Value = "La Araucanía"
$ws.Cells.Item(95, 4).Value = 44769
$ws.Cells.Item(95, 5).Value = 9
$ws.Cells.Item(95, 6).Value = 100114002
$ws.Cells.Item(95, 7).Value = "Camote"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 80
$ws.Cells.Item(95, 11).Value = 20000
$ws.Cells.Item(95, 12).Value = 20000
$ws.Cells.Item(95, 13).Value = 20000
$ws.Cells.Item(95, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(95, 15).Value = "Perú"
$ws.Cells.Item(95, 16).Value = 1000
$ws.Cells.Item(95, 17).Value = 20
$ws.Cells.Item(95, 18).Value = "Hortaliza"
